$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated input values in column G (Shared) ---
$ws.Range("G3").Value2 = 1750297
$ws.Range("G4").Value2 = 33867573
$ws.Range("G5").Value2 = 31353825
$ws.Range("G6").Value2 = 18968058

# --- New (blank, italic-styled) cells at G1 and H1 ---
$ws.Range("G1:H1").Font.Italic = $true

# --- Columns reset to the (new) default width instead of per-column overrides ---
$ws.StandardWidth = 11.85546875

# --- Cursor / selection moved to M9 ---
$ws.Range("M9").Select()
